$wb = $excel.ActiveWorkbook

# Both "展览" and "全部类型" sheets contain the same rows of convention data
# and need their "想去人数" (F column) counts bumped.
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 334
    $ws.Range("F4").Value = 179
    $ws.Range("F5").Value = 127
}
